$d = $word.ActiveDocument

# 1. ":ЗУ1" -> "69:32:0070101:ЗУ1" (3 occurrences)
$d.Content.Find.Execute(":ЗУ1", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "69:32:0070101:ЗУ1", 2)

# 2. "Метод спутниковых геодезических измерений (определений)" -> "None"
$d.Content.Find.Execute("Метод спутниковых геодезических измерений (определений)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "None", 2)

# 3. "ПО Trimble Business Center " -> "Mt = 0.1"
$d.Content.Find.Execute("ПО Trimble Business Center ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Mt = 0.1", 2)

# 4. "15" (exact) -> "667"
$d.Content.Find.Execute("15", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "667", 2)

# 5. "∆Р=2mt√p " -> "∆Р=3.5mt√p=3.5*0.1*√667=9"
$d.Content.Find.Execute("∆Р=2mt√p ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "∆Р=3.5mt√p=3.5*0.1*√667=9", 2)
